# Refresh market-price / profit figures pulled in by the scheduled runner.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) for the affected leve rows across the ALC/ARM/BSM/CRP/CUL/
# GSM/LTW sheets, matching the latest market snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 139.90909
$ws.Range("I11").Value = 139.90909
$ws.Range("K11").Value = 139.90909
$ws.Range("M11").Value = 0.09091000000000804
$ws.Range("H18").Value = 944.7857
$ws.Range("I18").Value = 957.5
$ws.Range("J18").Value = 868.5
$ws.Range("K18").Value = 957.5
$ws.Range("L18").Value = 868.5
$ws.Range("M18").Value = -673.5
$ws.Range("N18").Value = -1436.5
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H74").Value = 6000
$ws.Range("I74").Value = 2000
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 2000
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -1064
$ws.Range("N74").Value = -11872
$ws.Range("H77").Value = 6000
$ws.Range("I77").Value = 2000
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 10000
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -5320
$ws.Range("N77").Value = -59360
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H98").Value = 1043.7778
$ws.Range("I98").Value = 1142.1428
$ws.Range("J98").Value = 699.5
$ws.Range("K98").Value = 1142.1428
$ws.Range("L98").Value = 699.5
$ws.Range("M98").Value = 355.8571999999999
$ws.Range("N98").Value = -3695.5
$ws.Range("H122").Value = 1043.7778
$ws.Range("I122").Value = 1142.1428
$ws.Range("J122").Value = 699.5
$ws.Range("K122").Value = 3426.4284
$ws.Range("L122").Value = 2098.5
$ws.Range("M122").Value = -976.4284000000002
$ws.Range("N122").Value = -6998.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("N34").Value = 0
$ws.Range("H132").Value = 2473.6924
$ws.Range("I132").Value = 2367.1
$ws.Range("K132").Value = 7101.299999999999
$ws.Range("M132").Value = -4571.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 53269.875
$ws.Range("I107").Value = 69744.836
$ws.Range("K107").Value = 69744.836
$ws.Range("M107").Value = -67824.836
$ws.Range("H110").Value = 77000
$ws.Range("I110").Value = 10000
$ws.Range("J110").Value = 99333.336
$ws.Range("K110").Value = 10000
$ws.Range("L110").Value = 99333.336
$ws.Range("M110").Value = -5910
$ws.Range("N110").Value = -107513.336
$ws.Range("H130").Value = 30780
$ws.Range("J130").Value = 30780
$ws.Range("L130").Value = 30780
$ws.Range("N130").Value = -40820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 239500
$ws.Range("J9").Value = 239500
$ws.Range("L9").Value = 239500
$ws.Range("N9").Value = -239836
$ws.Range("H16").Value = 1207.25
$ws.Range("I16").Value = 1078
$ws.Range("J16").Value = 1595
$ws.Range("K16").Value = 1078
$ws.Range("L16").Value = 1595
$ws.Range("M16").Value = -791
$ws.Range("N16").Value = -2169
$ws.Range("H31").Value = 2667.88
$ws.Range("I31").Value = 1284.95
$ws.Range("K31").Value = 1284.95
$ws.Range("M31").Value = -989.95
$ws.Range("H34").Value = 2667.88
$ws.Range("I34").Value = 1284.95
$ws.Range("K34").Value = 1284.95
$ws.Range("M34").Value = -1082.95
$ws.Range("H58").Value = 1200
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1200
$ws.Range("K58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("M58").Value = 1200
$ws.Range("N58").Value = -1606
$ws.Range("H105").Value = 564.1667
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 750
$ws.Range("J107").Value = 700
$ws.Range("L107").Value = 700
$ws.Range("N107").Value = -4540
$ws.Range("H113").Value = 1207.25
$ws.Range("I113").Value = 1078
$ws.Range("J113").Value = 1595
$ws.Range("K113").Value = 1078
$ws.Range("L113").Value = 1595
$ws.Range("M113").Value = 1092
$ws.Range("N113").Value = -5935
$ws.Range("H136").Value = 1200
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").Value = 3600
$ws.Range("N136").Value = -8700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2639.8545
$ws.Range("J34").Value = 2723.415
$ws.Range("L34").Value = 8170.245
$ws.Range("N34").Value = -8338.244999999999
$ws.Range("H39").Value = 1768.8462
$ws.Range("J39").Value = 1768.8462
$ws.Range("L39").Value = 5306.5386
$ws.Range("N39").Value = -5894.5386
$ws.Range("H80").Value = 4915.5
$ws.Range("I80").Value = 4850
$ws.Range("K80").Value = 14550
$ws.Range("M80").Value = -13614
$ws.Range("H83").Value = 4915.5
$ws.Range("I83").Value = 4850
$ws.Range("K83").Value = 43650
$ws.Range("M83").Value = -38970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 37038052
$ws.Range("I107").Value = 157.25
$ws.Range("J107").Value = 66668370
$ws.Range("K107").Value = 157.25
$ws.Range("L107").Value = 66668370
$ws.Range("M107").Value = 1762.75
$ws.Range("N107").Value = -66672210
$ws.Range("H113").Value = 524.5
$ws.Range("I113").Value = 524.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 524.5
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1645.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 974.2273
$ws.Range("I22").Value = 786.0833
$ws.Range("K22").Value = 786.0833
$ws.Range("M22").Value = -491.0833
$ws.Range("H27").Value = 974.2273
$ws.Range("I27").Value = 786.0833
$ws.Range("K27").Value = 786.0833
$ws.Range("M27").Value = -679.0833
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10590
$ws.Range("H38").Value = 33033
$ws.Range("J38").Value = 33033
$ws.Range("L38").Value = 33033
$ws.Range("N38").Value = -33853
$ws.Range("H55").Value = 932
$ws.Range("I55").Value = 546.2
$ws.Range("J55").Value = 1575
$ws.Range("K55").Value = 546.2
$ws.Range("L55").Value = 1575
$ws.Range("M55").Value = -373.2
$ws.Range("N55").Value = -1921
$ws.Range("H61").Value = 4057.8333
$ws.Range("I61").Value = 3469.4
$ws.Range("K61").Value = 3469.4
$ws.Range("M61").Value = -3267.4
$ws.Range("H113").Value = 4057.8333
$ws.Range("I113").Value = 3469.4
$ws.Range("K113").Value = 3469.4
$ws.Range("M113").Value = -1299.4
$ws.Range("H136").Value = 1834.6666
$ws.Range("I136").Value = 304
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 912
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = 1638
$ws.Range("N136").Value = -12900
